# Auto update Excel log
# Appends new sensor/alert rows to the "ALERTS" and "PIR" sheets.

$wb = $excel.ActiveWorkbook

# ---- ALERTS sheet: append row 3 -----------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")

$alertsRow = 3
$alerts.Range("A$alertsRow").NumberFormat = "@"
$alerts.Range("A$alertsRow").Value = "2026-02-06"
$alerts.Range("A$alertsRow").ClearFormats()
$alerts.Range("B$alertsRow").Value = "09:42:30"
$alerts.Range("C$alertsRow").Value = "09:00"
$alerts.Range("D$alertsRow").Value = "Bathroom"
$alerts.Range("E$alertsRow").Value = "MINIMAL"
$alerts.Range("F$alertsRow").Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

# ---- PIR sheet: append rows 62-70 ---------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-02-06", "09:41:41", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:41:41", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:41:46", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:41:51", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:41:56", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:41:58", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:42:06", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:42:06", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:42:14", "09:00", "Bathroom", "No Motion",       "Inactive")
)

$startRow = 62
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $values = $pirRows[$i]

    $pir.Range("A$r").NumberFormat = "@"
    $pir.Range("A$r").Value = $values[0]
    $pir.Range("A$r").ClearFormats()

    $pir.Range("B$r").Value = $values[1]
    $pir.Range("C$r").Value = $values[2]
    $pir.Range("D$r").Value = $values[3]
    $pir.Range("E$r").Value = $values[4]
    $pir.Range("F$r").Value = $values[5]
}
